$d = $word.ActiveDocument

# --- 1. Title: "Project Phases Definition" -> "Project Phases "
#    plus the "_GoBack" bookmark, which moves here from further down
#    in the document (Word only keeps one bookmark per name, so
#    re-adding it here also removes the old one). ---

$prefix = "Project Phases "

$titleRange = $d.Content
$titleRange.Find.Execute("Project Phases Definition", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$splitPos = $titleRange.Start + $prefix.Length
$titleEnd = $titleRange.End

$bm = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bm)

$tail = $d.Range($splitPos, $titleEnd)
$tail.Delete()

# --- 2. Heading "Phase B - Learning the ropes"
#    -> "Phase B - " / "setting up base camp" as two separate runs. ---

$dash = [char]0x2013

$headingRange = $d.Content
$headingRange.Find.Execute("Learning the ropes", $true, $false, $false, $false, $false, $true, 1, $false, "setting up base camp", 2)

$prefixRange = $d.Content
$prefixRange.Find.Execute("Phase B " + $dash + " ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Force a run boundary right after "Phase B - " using the same
# add-then-delete bookmark trick as above, leaving no bookmark behind.
$splitPos2 = $prefixRange.End
$bm2 = $d.Range($splitPos2, $splitPos2)
$d.Bookmarks.Add("TEMP_RUN_SPLIT", $bm2)
$d.Bookmarks("TEMP_RUN_SPLIT").Delete()
